$wb = $excel.ActiveWorkbook

# Overview sheet: update "Latest Handoff Date" for af41e3fc-efcf-45dd-9817-8262d2f016a0.md (row 5)
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-13-17 06:13:42"

# zh-cn sheet: update "Latest Handoff Datetime" for the same file (row 5)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-17 06:13:34"

# de-de sheet: update "Latest Handoff Datetime" for the same file (row 5)
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-17 06:13:42"
